$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.880.15"
$ws.Range("E2").Value = "  +1.93%  "
$ws.Range("D3").Value = "1.708.19"
$ws.Range("E3").Value = "  +1.73%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'313.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.99%  "
$ws.Range("D6").Value = "'0.9988"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").Value = "'0.3745"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.92%  "
$ws.Range("D8").Value = "'49.44"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.54%  "
$ws.Range("D9").Value = "'0.3439"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("D10").Value = "'1.225"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.63%  "
$ws.Range("D11").Value = "'0.07546"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.89%  "
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("D13").Value = "'21.29"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.23%  "
$ws.Range("D14").Value = "'6.321"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.46%  "
$ws.Range("D15").Value = "'7.080"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.97%  "
$ws.Range("D16").Value = "1.708.73"
$ws.Range("E16").Value = "  +1.70%  "
$ws.Range("D17").Value = "'0.00001131"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.22%  "
$ws.Range("D18").Value = "'0.06732"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("D19").Value = "'0.9984"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("D20").Value = "'84.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.36%  "
$ws.Range("D21").Value = "'17.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.26%  "
$ws.Range("D22").Value = "'6.395"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.50%  "
$ws.Range("D23").Value = "'13.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.53%  "
$ws.Range("D24").Value = "24.902.55"
$ws.Range("E24").Value = "  +2.21%  "
$ws.Range("D25").Value = "'2.449"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.45%  "
$ws.Range("D26").Value = "'2.800"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.46%  "
$ws.Range("D27").Value = "'20.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.64%  "
$ws.Range("E28").Value = "  -2.64%  "
$ws.Range("D29").Value = "'132.91"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.21%  "
$ws.Range("D30").Value = "1.896.78"
$ws.Range("E30").Value = "  +1.68%  "
$ws.Range("D31").Value = "'1.248"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +28.22%  "
$ws.Range("D32").Value = "'6.828"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.22%  "
$ws.Range("D33").Value = "'4.226"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.04%  "
$ws.Range("D34").Value = "'13.96"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +13.24%  "
$ws.Range("B35").Value = "WEMIXTOKEN"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "'1.780"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.82%  "
$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D36").Value = "'0.08802"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.24%  "
$ws.Range("D37").Value = "'5.640"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.61%  "
$ws.Range("D38").Value = "'0.06656"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.79%  "
$ws.Range("D39").Value = "'9.181"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.34%  "
$ws.Range("D40").Value = "'0.02412"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.88%  "
$ws.Range("D41").Value = "'0.2234"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.54%  "
$ws.Range("D42").Value = "'1.277"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.02%  "
$ws.Range("D43").Value = "'0.6466"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.15%  "
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").Value = "'13.90"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.38%  "
$ws.Range("D46").Value = "'0.6162"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.96%  "
$ws.Range("E47").Value = "  +1.65%  "
$ws.Range("D48").Value = "'2.129"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.23%  "
$ws.Range("D49").Value = "'129.64"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.90%  "
$ws.Range("D50").Value = "'0.07320"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.93%  "
$ws.Range("D51").Value = "'80.13"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.78%  "
